# Upload Y4_B2526_Excuses.xlsx via attendance app
# The attendance app re-synced the log: the previous second excuse entry
# (row 3, student 211242) is gone, and the remaining entry (row 2) now
# reflects a newer log line (different student id / log date).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Excuses")

# --- Update row 2 in place -------------------------------------------------
# Student ID (A2): keep it text ("@" number format) so "211210" isn't
# stored as a numeric value, matching the original inline-string cell.
$a2 = $ws.Range("A2")
$a2.NumberFormat = "@"
$a2.Value = "211210"

# Log Date (C2)
$ws.Range("C2").Value = "19/10/2025"

# Re-apply row 2's original formatting to A2 (copying it off a sibling
# cell in the same row) so the forced text format above doesn't leave
# behind a stray/duplicate cell style.
$ws.Range("B2").Copy()
$a2.PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# --- Remove row 3 entirely --------------------------------------------------
$ws.Rows.Item(3).Delete()
